$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 247; existing rows 247-298 shift down to 248-299,
# inheriting formatting (including the date number format on column D) from
# the surrounding rows, matching Excel's default Insert-row behaviour.
$ws.Rows.Item(247).Insert()

# Populate the newly inserted row 247 with the new weekly record.
$ws.Cells.Item(247, 1).Value = 4
$ws.Cells.Item(247, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(247, 3).Value = "Los Lagos"
$ws.Cells.Item(247, 4).Value = 44785
$ws.Cells.Item(247, 5).Value = 10
$ws.Cells.Item(247, 6).Value = 100112003
$ws.Cells.Item(247, 7).Value = "Ajo"
$ws.Cells.Item(247, 8).Value = "Chino"
$ws.Cells.Item(247, 9).Value = "Primera"
$ws.Cells.Item(247, 10).Value = 200
$ws.Cells.Item(247, 11).Value = 28000
$ws.Cells.Item(247, 12).Value = 29000
$ws.Cells.Item(247, 13).Value = 28500
$ws.Cells.Item(247, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(247, 15).Value = "China"
$ws.Cells.Item(247, 16).Value = 2850
$ws.Cells.Item(247, 17).Value = 10
$ws.Cells.Item(247, 18).Value = "Hortaliza"
